$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new "2023" data column (S) -----------------------------------
# Copy the formatting from column R (the previous last year column) into S,
# then fill in the new values for the rows that have data.

$ws.Range("R1:R15").Copy() | Out-Null
$ws.Range("S1:S15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 33.799999999999997
$ws.Range("S5").Value = 6.6
$ws.Range("S6").Value = 442
$ws.Range("S7").Value = 293
$ws.Range("S8").Value = 538
$ws.Range("S9").Value = 1.3
$ws.Range("S10").Value = 1.9
$ws.Range("S11").Value = 5.3
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 29
$ws.Range("S14").Value = 27.8

# --- Extend the trailing blank column (T -> U) -----------------------------
# The sheet always keeps one blank formatted column after the data; shift it
# from T to the new column U.

$ws.Range("T2:T15").Copy() | Out-Null
$ws.Range("U2:U15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Merged title cell now spans through column S --------------------------
$ws.Range("A1:R1").UnMerge() | Out-Null
$ws.Range("A1:S1").Merge() | Out-Null

# --- Selection moves to the new column --------------------------------------
$ws.Range("S3:S14").Select() | Out-Null

# --- Window size change recorded alongside the edit -------------------------
$excel.ActiveWindow.Width = 16605
$excel.ActiveWindow.Height = 8610
